$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.162.87"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "1.905.75"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.99"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4607"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3890"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07882"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9904"
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.96"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").Value = "1.909.88"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.769"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.044"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07041"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.20"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009954"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.07"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "29.161.41"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.320"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.22"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.47"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.909"
$ws.Range("E27").Value = "  -2.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "118.78"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.878"
$ws.Range("E29").Value = "  -6.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09355"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8955"
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.231"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.321"
$ws.Range("E33").Value = "  -2.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.164"
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05797"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.171"
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02087"
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.002"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5715"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.675"
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1805"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.717"
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.96"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5362"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.172"
$ws.Range("E45").Value = "  -4.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07023"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.550"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.25"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.2947"
$ws.Range("E50").Value = "  +0.34%  "

# Row 51: WEMIXToken -> Aave
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.30"
$ws.Range("E51").Value = "  -0.55%  "
